$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-13 (header row 1 stays the same)
# Columns: A, B (dataset name), C (execution_time-ish), D, E (method), F (n_tries), G (search_space)
$data = @(
    @(0,  "10_features_0.01_error.csv.csv",  7854.974847266059, 2.539521217346191, "bipop", 30, 0.1),
    @(1,  "5_features_0.01_error.csv.csv",   7384.869514808409, 2.529790163040161, "bipop", 30, 0.1),
    @(2,  "10_features_0.5_error.csv.csv",   9075.575259189287, 2.983652114868164, "bipop", 30, 0.1),
    @(3,  "10_features_0.1_error.csv.csv",   8225.308119585436, 3.23308801651001,  "bipop", 30, 0.1),
    @(4,  "100_features_0.1_error.csv.csv",  7997.539347149157, 3.222321748733521, "bipop", 30, 0.1),
    @(5,  "100_features_0.5_error.csv.csv",  9014.655831042939, 2.663567066192627, "bipop", 30, 0.1),
    @(6,  "5_features_0.5_error.csv.csv",    9310.509594784353, 2.627505779266357, "bipop", 30, 0.1),
    @(7,  "100_features_0.01_error.csv.csv", 8046.293624039642, 2.640698909759521, "bipop", 30, 0.1),
    @(8,  "5_features_0.1_error.csv.csv",    8224.794069499589, 2.661087989807129, "bipop", 30, 0.1),
    @(9,  "2_features_0.01_error.csv.csv",   7309.910203919307, 3.539085865020752, "bipop", 30, 0.1),
    @(10, "2_features_0.1_error.csv.csv",    8609.001889490926, 2.680201768875122, "bipop", 30, 0.1),
    @(11, "2_features_0.5_error.csv.csv",    9812.223000884836, 2.807796955108643, "bipop", 30, 0.1)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $row++
}

# Remove the now-obsolete rows 14, 15, 16 (old sheet had data through row 16)
$ws.Range("A14:G16").Delete()
